# Fixed update to excel issue
# - Rename "Requested quantity" headers to metric-specific names
# - Add a new "PO Forecast" sheet with forecast data

$wb = $excel.ActiveWorkbook

# --- Weekly Quantity sheet: rename header ---
$wsWeekly = $wb.Worksheets.Item(1)
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# --- Monthly Trend sheet: rename header ---
$wsMonthly = $wb.Worksheets.Item(2)
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy header formatting (bold/border style) from an existing header cell
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Copy date-column formatting from the existing date column
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A12").PasteSpecial(-4122)

# Header labels
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Data rows
$wsForecast.Cells.Item(2, 1).Value = 45578.99999999999
$wsForecast.Cells.Item(2, 2).Value = 654
$wsForecast.Cells.Item(2, 3).Value = 501.4200453919026
$wsForecast.Cells.Item(2, 4).Value = 810.5514582527992

$wsForecast.Cells.Item(3, 1).Value = 45592.99999999999
$wsForecast.Cells.Item(3, 2).Value = 196
$wsForecast.Cells.Item(3, 3).Value = 57.23901610424889
$wsForecast.Cells.Item(3, 4).Value = 338.648632713941

$wsForecast.Cells.Item(4, 1).Value = 45599.99999999999
$wsForecast.Cells.Item(4, 2).Value = 0
$wsForecast.Cells.Item(4, 3).Value = -184.4082983346196
$wsForecast.Cells.Item(4, 4).Value = 109.8278708457326

$wsForecast.Cells.Item(5, 1).Value = 45606.99999999999
$wsForecast.Cells.Item(5, 2).Value = 0
$wsForecast.Cells.Item(5, 3).Value = -411.6510667299753
$wsForecast.Cells.Item(5, 4).Value = -118.4878096971054

$wsForecast.Cells.Item(6, 1).Value = 45613.99999999999
$wsForecast.Cells.Item(6, 2).Value = 0
$wsForecast.Cells.Item(6, 3).Value = -629.745396146613
$wsForecast.Cells.Item(6, 4).Value = -347.369769862345

$wsForecast.Cells.Item(7, 1).Value = 45620.99999999999
$wsForecast.Cells.Item(7, 2).Value = 0
$wsForecast.Cells.Item(7, 3).Value = -862.6112439083638
$wsForecast.Cells.Item(7, 4).Value = -573.7356368138794

$wsForecast.Cells.Item(8, 1).Value = 45627.99999999999
$wsForecast.Cells.Item(8, 2).Value = 0
$wsForecast.Cells.Item(8, 3).Value = -1099.309379065477
$wsForecast.Cells.Item(8, 4).Value = -812.8197957718797

$wsForecast.Cells.Item(9, 1).Value = 45634.99999999999
$wsForecast.Cells.Item(9, 2).Value = 0
$wsForecast.Cells.Item(9, 3).Value = -1319.003702599054
$wsForecast.Cells.Item(9, 4).Value = -1030.802717229441

$wsForecast.Cells.Item(10, 1).Value = 45641.99999999999
$wsForecast.Cells.Item(10, 2).Value = 0
$wsForecast.Cells.Item(10, 3).Value = -1555.661274064956
$wsForecast.Cells.Item(10, 4).Value = -1269.136341623547

$wsForecast.Cells.Item(11, 1).Value = 45648.99999999999
$wsForecast.Cells.Item(11, 2).Value = 0
$wsForecast.Cells.Item(11, 3).Value = -1792.506212419167
$wsForecast.Cells.Item(11, 4).Value = -1496.591210447476

$wsForecast.Cells.Item(12, 1).Value = 45655.99999999999
$wsForecast.Cells.Item(12, 2).Value = 0
$wsForecast.Cells.Item(12, 3).Value = -2012.562560109675
$wsForecast.Cells.Item(12, 4).Value = -1726.142954676881

# Put selection back to A1 to match default sheet view state
$wsForecast.Range("A1").Select()
